# Add a "Priority" column to the Product Backlog sheet.
# This inserts a new column A (shifting the existing Features/Adjust
# Factor/Adjust Estimate/Initial Estimates columns one to the right),
# labels the new header "Priority" in A2 (bold, matching the other
# header cells), and numbers each backlog item 1-16 down column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift B:E (old A:D) one column to the right, carrying values, styles,
# formulas and the title merge cell along with it.
$ws.Columns("A").Insert()

# New header cell for the inserted column.
$ws.Range("A2").Value = "Priority"
$ws.Range("A2").Font.Bold = $true

# Sequential priority numbers for each of the 16 backlog rows (3-18).
for ($row = 3; $row -le 18; $row++) {
    $ws.Cells.Item($row, 1).Value = $row - 2
}

# Re-apply the "Adjust Estimate" formula as a single range assignment so
# the fill keeps using one shared formula across E8:E18, same as before
# the column insert.
$ws.Range("E8:E18").Formula = "=C8+(C8*D8)"

$ws.Range("A2").Select()
